$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startSerial = 44432
$startRow = 358

$data = @(
    @(1, 6, 244.798041615667),
    @(0, 6, 244.798041615667),
    @(1, 5, 203.9983680130559),
    @(0, 4, 163.1986944104447),
    @(1, 5, 203.9983680130559),
    @(0, 5, 203.9983680130559),
    @(0, 3, 122.3990208078335),
    @(0, 2, 81.59934720522236),
    @(0, 2, 81.59934720522236)
)

# First, replicate the formatting used by the last pre-existing row's column A
# cell (style s="2") onto the new column-A cells, before writing their values.
$srcCell = $ws.Cells.Item($startRow - 1, 1)
$destRange = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($startRow + $data.Length - 1, 1))
$srcCell.Copy($destRange)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $serial = $startSerial + $i
    $ws.Cells.Item($row, 1).Value = $serial
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
}
